$wb = $excel.ActiveWorkbook

$wsAstronauta = $wb.Worksheets.Item("Astronauta")
$wsMago       = $wb.Worksheets.Item("Mago")
$wsNinja      = $wb.Worksheets.Item("Ninja")

# Add new attitude scores in column C for the "Ninja" sheet
$wsNinja.Range("C2").Value = 1
$wsNinja.Range("C3").Value = 1
$wsNinja.Range("C4").Value = 1
$wsNinja.Range("C5").Value = 0
$wsNinja.Range("C6").Value = 1
$wsNinja.Range("C7").Value = 1

# Update the selection on the "Mago" sheet (no longer the active tab)
$wsMago.Range("B7").Select()

# Update the selection on the "Ninja" sheet
$wsNinja.Range("C8").Select()

# Make "Astronauta" the active sheet/tab and set its selection
$wsAstronauta.Activate()
$wsAstronauta.Range("C7").Select()
